# Update pollutant names to lowercase (shared strings) and refresh the
# recomputed incidence values (column D) to their new precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: pollutant names -> lowercase
$ws.Range("B3").Value = "ethyl benzene"
$ws.Range("B4").Value = "ethyl benzene"
$ws.Range("B5").Value = "toluene"
$ws.Range("B6").Value = "xylenes (mixed)"

# Column D: updated Incidence values (precision refresh)
$newIncidence = [double]"1.80935905528e-06"
$ws.Range("D2").Value = $newIncidence
$ws.Range("D3").Value = $newIncidence
$ws.Range("D4").Value = $newIncidence
$ws.Range("D7").Value = $newIncidence
